$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3798.889
$ws.Range("I116").Value = 4257.143
$ws.Range("J116").Value = 2195
$ws.Range("K116").Value = 4257.143
$ws.Range("L116").Value = 2195
$ws.Range("M116").Value = -815.143
$ws.Range("N116").Value = -9079
$ws.Range("H137").Value = 37039308
$ws.Range("I137").Value = 55557320
$ws.Range("J137").Value = 3287.3333
$ws.Range("K137").Value = 166671960
$ws.Range("L137").Value = 9861.999899999999
$ws.Range("M137").Value = -166669410
$ws.Range("N137").Value = -14961.9999
$ws.Range("H138").Value = 10396476
$ws.Range("I138").Value = 5956829
$ws.Range("J138").Value = 11635447
$ws.Range("K138").Value = 17870487
$ws.Range("L138").Value = 34906341
$ws.Range("M138").Value = -17865347
$ws.Range("N138").Value = -34916621

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H32").Value = 60748.77
$ws.Range("I32").Value = 13157
$ws.Range("K32").Value = 13157
$ws.Range("M32").Value = -12870
$ws.Range("H61").Value = 1956.0555
$ws.Range("I61").Value = 1705.0938
$ws.Range("K61").Value = 1705.0938
$ws.Range("M61").Value = -1493.0938
$ws.Range("H74").Value = 6116.6665
$ws.Range("I74").Value = 1274.1428
$ws.Range("J74").Value = 11331.692
$ws.Range("K74").Value = 1274.1428
$ws.Range("L74").Value = 11331.692
$ws.Range("M74").Value = -400.1428000000001
$ws.Range("N74").Value = -13079.692
$ws.Range("H77").Value = 6116.6665
$ws.Range("I77").Value = 1274.1428
$ws.Range("J77").Value = 11331.692
$ws.Range("K77").Value = 6370.714
$ws.Range("L77").Value = 56658.45999999999
$ws.Range("M77").Value = -2002.714
$ws.Range("N77").Value = -65394.45999999999
$ws.Range("H136").Value = 1956.0555
$ws.Range("I136").Value = 1705.0938
$ws.Range("K136").Value = 5115.2814
$ws.Range("M136").Value = -2565.2814

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 2016
$ws.Range("I29").Value = 2016
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 2016
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -1727
$ws.Range("N29").ClearContents()
$ws.Range("H36").Value = 2500
$ws.Range("I36").Value = 2500
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2500
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1966
$ws.Range("N36").ClearContents()
$ws.Range("H75").Value = 133719.5
$ws.Range("I75").Value = 11547.333
$ws.Range("K75").Value = 11547.333
$ws.Range("M75").Value = -10611.333
$ws.Range("H78").Value = 133719.5
$ws.Range("I78").Value = 11547.333
$ws.Range("K78").Value = 34641.999
$ws.Range("M78").Value = -29961.999
$ws.Range("H105").Value = 3066.75
$ws.Range("I105").Value = 2861.0344
$ws.Range("K105").Value = 2861.0344
$ws.Range("M105").Value = -1114.0344

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4259.3237
$ws.Range("I31").Value = 2253.7083
$ws.Range("J31").Value = 9072.799999999999
$ws.Range("K31").Value = 2253.7083
$ws.Range("L31").Value = 9072.799999999999
$ws.Range("M31").Value = -1958.7083
$ws.Range("N31").Value = -9662.799999999999
$ws.Range("H34").Value = 4259.3237
$ws.Range("I34").Value = 2253.7083
$ws.Range("J34").Value = 9072.799999999999
$ws.Range("K34").Value = 2253.7083
$ws.Range("L34").Value = 9072.799999999999
$ws.Range("M34").Value = -2051.7083
$ws.Range("N34").Value = -9476.799999999999
$ws.Range("H132").Value = 5379391.5
$ws.Range("I132").Value = 8774016
$ws.Range("J132").Value = 4570.25
$ws.Range("K132").Value = 26322048
$ws.Range("L132").Value = 13710.75
$ws.Range("M132").Value = -26319518
$ws.Range("N132").Value = -18770.75

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1605.5135
$ws.Range("I5").Value = 669.10345
$ws.Range("K5").Value = 2007.31035
$ws.Range("M5").Value = -1895.31035
$ws.Range("H10").Value = 3100
$ws.Range("I10").Value = 200
$ws.Range("K10").Value = 600
$ws.Range("M10").Value = -461
$ws.Range("H94").Value = 5000
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H105").Value = 5333.3335
$ws.Range("J105").Value = 5333.3335
$ws.Range("L105").Value = 16000.0005
$ws.Range("N105").Value = -21242.0005
$ws.Range("H113").Value = 5682380.5
$ws.Range("I113").Value = 550
$ws.Range("J113").Value = 7576324
$ws.Range("K113").Value = 1650
$ws.Range("L113").Value = 22728972
$ws.Range("M113").Value = 520
$ws.Range("N113").Value = -22733312
$ws.Range("H124").Value = 333334370
$ws.Range("I124").Value = 1530
$ws.Range("K124").Value = 4590
$ws.Range("M124").Value = 320
$ws.Range("H125").Value = 2933.1035
$ws.Range("I125").Value = 2030
$ws.Range("K125").Value = 6090
$ws.Range("M125").Value = -1170
$ws.Range("H131").Value = 3624544.5
$ws.Range("J131").Value = 3789264.8
$ws.Range("L131").Value = 11367794.4
$ws.Range("N131").Value = -11377874.4
$ws.Range("H135").Value = 1605.5135
$ws.Range("I135").Value = 669.10345
$ws.Range("K135").Value = 6021.931049999999
$ws.Range("M135").Value = -3486.931049999999

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 2000000
$ws.Range("I7").Value = 2000000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2000000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1999888
$ws.Range("N7").ClearContents()
$ws.Range("H8").Value = 2000000
$ws.Range("I8").Value = 2000000
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 2000000
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -1999861
$ws.Range("N8").ClearContents()
$ws.Range("H43").Value = 25400
$ws.Range("J43").Value = 19250
$ws.Range("L43").Value = 19250
$ws.Range("N43").Value = -19552
$ws.Range("H46").Value = 18213.615
$ws.Range("J46").Value = 19500
$ws.Range("L46").Value = 19500
$ws.Range("N46").Value = -19812
$ws.Range("H64").Value = 10000
$ws.Range("I64").Value = 10000
$ws.Range("K64").Value = 10000
$ws.Range("M64").Value = -9752
$ws.Range("H67").Value = 10000
$ws.Range("I67").Value = 10000
$ws.Range("K67").Value = 10000
$ws.Range("M67").Value = -9142
$ws.Range("H80").Value = 68184744
$ws.Range("I80").Value = 3242.8572
$ws.Range("J80").Value = 187502370
$ws.Range("K80").Value = 3242.8572
$ws.Range("L80").Value = 187502370
$ws.Range("M80").Value = -2244.8572
$ws.Range("N80").Value = -187504366
$ws.Range("H83").Value = 68184744
$ws.Range("I83").Value = 3242.8572
$ws.Range("J83").Value = 187502370
$ws.Range("K83").Value = 16214.286
$ws.Range("L83").Value = 937511850
$ws.Range("M83").Value = -11222.286
$ws.Range("N83").Value = -937521834
$ws.Range("H97").Value = 1985.5
$ws.Range("I97").Value = 1999.6666
$ws.Range("J97").Value = 1960
$ws.Range("K97").Value = 1999.6666
$ws.Range("L97").Value = 1960
$ws.Range("M97").Value = -1503.6666
$ws.Range("N97").Value = -2952
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6030.4736
$ws.Range("I136").Value = 2448.5
$ws.Range("J136").Value = 16060
$ws.Range("K136").Value = 7345.5
$ws.Range("L136").Value = 48180
$ws.Range("M136").Value = -4795.5
$ws.Range("N136").Value = -53280

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 40531.84
$ws.Range("I126").Value = 59260.766
$ws.Range("K126").Value = 177782.298
$ws.Range("M126").Value = -175312.298
$ws.Range("H136").Value = 2672.8975
$ws.Range("I136").Value = 726.375
$ws.Range("J136").Value = 11571.286
$ws.Range("K136").Value = 2179.125
$ws.Range("L136").Value = 34713.858
$ws.Range("M136").Value = 370.875
$ws.Range("N136").Value = -39813.858
